# Auto-generated edit script: updates cryptos price/volume data
# to match the "Tue Mar 28 22:51:31 UTC 2023" GitHub Actions refresh commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.239.08"
$ws.Range("E2").Value = "  +0.35%  "

$ws.Range("D3").Value = "1.772.15"
$ws.Range("E3").Value = "  +3.53%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'312.83"
$ws.Range("E5").Value = "  +1.29%  "

$ws.Range("E6").Value = "  +0.05%  "

$ws.Range("D7").Value = "'0.5221"
$ws.Range("E7").Value = "  +9.60%  "

$ws.Range("D8").Value = "'0.3670"
$ws.Range("E8").Value = "  +6.72%  "

$ws.Range("D9").Value = "'42.70"
$ws.Range("E9").Value = "  +1.34%  "

$ws.Range("D10").Value = "'0.07351"
$ws.Range("E10").Value = "  +0.99%  "

$ws.Range("E11").Value = "  +4.45%  "

$ws.Range("D12").Value = "'1.002"
$ws.Range("E12").Value = "  +0.10%  "

$ws.Range("D13").Value = "'20.43"
$ws.Range("E13").Value = "  +3.07%  "

$ws.Range("D14").Value = "'6.055"
$ws.Range("E14").Value = "  +3.48%  "

$ws.Range("D15").Value = "1.771.57"
$ws.Range("E15").Value = "  +3.65%  "

$ws.Range("D16").Value = "'6.925"
$ws.Range("E16").Value = "  +1.32%  "

$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("E18").Value = "  +0.15%  "

$ws.Range("E19").Value = "  +1.24%  "

$ws.Range("E20").Value = "  +0.05%  "

$ws.Range("E21").Value = "  +1.32%  "

$ws.Range("D22").Value = "'5.805"

$ws.Range("D23").Value = "27.273.17"
$ws.Range("E23").Value = "  +0.36%  "

$ws.Range("E24").Value = "  +4.08%  "

$ws.Range("D25").Value = "'2.132"
$ws.Range("E25").Value = "  +1.98%  "

$ws.Range("D26").Value = "'154.98"
$ws.Range("E26").Value = "  +1.38%  "

$ws.Range("E27").Value = "  +2.22%  "

$ws.Range("D28").Value = "1.973.50"
$ws.Range("E28").Value = "  +3.61%  "

$ws.Range("D29").Value = "'2.315"
$ws.Range("E29").Value = "  +11.28%  "

$ws.Range("D30").Value = "'120.98"
$ws.Range("E30").Value = "  +0.79%  "

$ws.Range("D31").Value = "'1.054"
$ws.Range("E31").Value = "  +4.04%  "

$ws.Range("D32").Value = "'0.09765"
$ws.Range("E32").Value = "  +5.42%  "

$ws.Range("D33").Value = "'5.547"
$ws.Range("E33").Value = "  +4.83%  "

$ws.Range("D34").Value = "'3.623"
$ws.Range("E34").Value = "  +0.90%  "

$ws.Range("D35").Value = "'0.02229"
$ws.Range("E35").Value = "  +1.65%  "

$ws.Range("D36").Value = "'0.05946"
$ws.Range("E36").Value = "  +0.90%  "

$ws.Range("D37").Value = "'11.18"
$ws.Range("E37").Value = "  +1.23%  "

$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'4.817"
$ws.Range("E38").Value = "  +1.40%  "

$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.6115"
$ws.Range("E39").Value = "  +3.32%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2014"
$ws.Range("E40").Value = "  +0.22%  "

$ws.Range("D41").Value = "'1.430"
$ws.Range("E41").Value = "  +1.33%  "

$ws.Range("D42").Value = "'8.045"
$ws.Range("E42").Value = "  +7.37%  "

$ws.Range("D43").Value = "'1.134"
$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("D44").Value = "'13.13"
$ws.Range("E44").Value = "  +3.97%  "

$ws.Range("D45").Value = "'0.5750"
$ws.Range("E45").Value = "  +2.40%  "

$ws.Range("D46").Value = "'3.623"
$ws.Range("E46").Value = "  +1.62%  "

$ws.Range("D47").Value = "'121.08"
$ws.Range("E47").Value = "  +2.38%  "

$ws.Range("D48").Value = "'1.873"
$ws.Range("E48").Value = "  +1.97%  "

$ws.Range("D49").Value = "'1.112"
$ws.Range("E49").Value = "  +2.47%  "

$ws.Range("D50").Value = "'0.06700"
$ws.Range("E50").Value = "  +1.01%  "

$ws.Range("D51").Value = "'1.002"
$ws.Range("E51").Value = "  +0.10%  "
